# Apply updated cryptocurrency price/volume figures (columns D and E)
# for rows 2-51, matching the scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.346.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "'2.485.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.07%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'565.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").Value = "'163.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.59%  "

$ws.Range("D9").Value = "'2.483.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("D10").Value = "'0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("D13").Value = "'4.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").Value = "'2.942.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("D15").Value = "'69.235.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("D16").Value = "'0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").Value = "'24.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.35%  "

$ws.Range("D18").Value = "'2.482.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.80%  "

$ws.Range("D19").Value = "'11.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").Value = "'7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.77%  "

$ws.Range("D21").Value = "'343.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.59%  "

$ws.Range("D22").Value = "'3.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").Value = "'1.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.00%  "

$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").Value = "'69.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").Value = "'3.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.01%  "

$ws.Range("D27").Value = "'2.613.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "

$ws.Range("D28").Value = "'8.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").Value = "'0.0₃0868"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.02%  "

$ws.Range("D31").Value = "'7.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.65%  "

$ws.Range("D32").Value = "'440.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("E33").Value = "  -6.28%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").Value = "'1.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").Value = "'155.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D37").Value = "'0.112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.08%  "

$ws.Range("D38").Value = "'19.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'18.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.89%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "'0.313"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("D42").Value = "'4.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.33%  "

$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "'2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.52%  "

$ws.Range("E45").Value = "  -7.79%  "

$ws.Range("D46").Value = "'137.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.28%  "

$ws.Range("D47").Value = "'3.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").Value = "'0.509"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").Value = "'0.0726"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").Value = "'0.569"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").Value = "'0.0919"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.11%  "
